$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fix existing rows 106/107: Spin Version 3 -> 4 ---
$ws.Cells.Item(106, 6).Value = 4
$ws.Cells.Item(107, 6).Value = 4

# --- Append new test case rows 108-129 ---
$ws.Cells.Item(108, 2).Value = "AR000"
$ws.Cells.Item(108, 3).Value = "Passed"
$ws.Cells.Item(108, 4).Value = "DFDL-17-007R"
$ws.Cells.Item(108, 5).Value = "High"
$ws.Cells.Item(108, 6).Value = 4
$ws.Cells.Item(108, 7).Value = "byte"
$ws.Cells.Item(108, 8).Value = "DFDL-166"

$ws.Cells.Item(109, 2).Value = "AQ000"
$ws.Cells.Item(109, 3).Value = "Passed"
$ws.Cells.Item(109, 4).Value = "DFDL-17-007R"
$ws.Cells.Item(109, 5).Value = "High"
$ws.Cells.Item(109, 6).Value = 4
$ws.Cells.Item(109, 7).Value = "text"
$ws.Cells.Item(109, 8).Value = "DFDL-165"

$ws.Cells.Item(110, 2).Value = "AA000"
$ws.Cells.Item(110, 3).Value = "Passed"
$ws.Cells.Item(110, 4).Value = "DFDL-17-007R"
$ws.Cells.Item(110, 5).Value = "High"
$ws.Cells.Item(110, 6).Value = 4
$ws.Cells.Item(110, 7).Value = "byte"
$ws.Cells.Item(110, 8).Value = "DFDL-148"

$ws.Cells.Item(111, 2).Value = "inputValueCalcErrorDiagnostic1"
$ws.Cells.Item(111, 3).Value = "Passed"
$ws.Cells.Item(111, 4).Value = "DFDL-17-007R"
$ws.Cells.Item(111, 5).Value = "High"
$ws.Cells.Item(111, 6).Value = 4
$ws.Cells.Item(111, 7).Value = "n/a"
$ws.Cells.Item(111, 8).Value = "DFDL-148"
$ws.Cells.Item(111, 9).Value = "Yes"

$ws.Cells.Item(112, 2).Value = "inputValueCalcErrorDiagnostic2"
$ws.Cells.Item(112, 3).Value = "Passed"
$ws.Cells.Item(112, 4).Value = "DFDL-17-007R"
$ws.Cells.Item(112, 5).Value = "High"
$ws.Cells.Item(112, 6).Value = 4
$ws.Cells.Item(112, 7).Value = "n/a"
$ws.Cells.Item(112, 8).Value = "DFDL-148"
$ws.Cells.Item(112, 9).Value = "Yes"

$ws.Cells.Item(113, 2).Value = "DelimProp_05"
$ws.Cells.Item(113, 3).Value = "Passed"
$ws.Cells.Item(113, 4).Value = "DFDL-12-033R"
$ws.Cells.Item(113, 5).Value = "High"
$ws.Cells.Item(113, 6).Value = 4
$ws.Cells.Item(113, 7).Value = "text"
$ws.Cells.Item(113, 8).Value = "DFDL-271"

$ws.Cells.Item(114, 2).Value = "AJ000"
$ws.Cells.Item(114, 3).Value = "Passed"
$ws.Cells.Item(114, 4).Value = "DFDL-5-015R"
$ws.Cells.Item(114, 5).Value = "High"
$ws.Cells.Item(114, 6).Value = 4
$ws.Cells.Item(114, 7).Value = "byte"
$ws.Cells.Item(114, 8).Value = "DFDL-159"

$ws.Cells.Item(115, 2).Value = "AJ001"
$ws.Cells.Item(115, 3).Value = "Passed"
$ws.Cells.Item(115, 4).Value = "DFDL-5-015R"
$ws.Cells.Item(115, 5).Value = "High"
$ws.Cells.Item(115, 6).Value = 4
$ws.Cells.Item(115, 7).Value = "byte"
$ws.Cells.Item(115, 8).Value = "DFDL-159"

$ws.Cells.Item(116, 2).Value = "text_02"
$ws.Cells.Item(116, 3).Value = "Passed"
$ws.Cells.Item(116, 4).Value = "DFDL-6-046R"
$ws.Cells.Item(116, 5).Value = "High"
$ws.Cells.Item(116, 6).Value = 3
$ws.Cells.Item(116, 7).Value = "text"
$ws.Cells.Item(116, 8).Value = "DFDL-199"

$ws.Cells.Item(117, 2).Value = "text_01"
$ws.Cells.Item(117, 3).Value = "Passed"
$ws.Cells.Item(117, 4).Value = "DFDL-13-235R"
$ws.Cells.Item(117, 5).Value = "High"
$ws.Cells.Item(117, 6).Value = 3
$ws.Cells.Item(117, 7).Value = "text"
$ws.Cells.Item(117, 8).Value = "DFDL-199"

$ws.Cells.Item(118, 2).Value = "text_03"
$ws.Cells.Item(118, 3).Value = "Passed"
$ws.Cells.Item(118, 4).Value = "DFDL-13-235R"
$ws.Cells.Item(118, 5).Value = "High"
$ws.Cells.Item(118, 6).Value = 3
$ws.Cells.Item(118, 7).Value = "text"
$ws.Cells.Item(118, 8).Value = "DFDL-199"

$ws.Cells.Item(119, 2).Value = "text_04"
$ws.Cells.Item(119, 3).Value = "Passed"
$ws.Cells.Item(119, 4).Value = "DFDL-13-235R"
$ws.Cells.Item(119, 5).Value = "High"
$ws.Cells.Item(119, 6).Value = 3
$ws.Cells.Item(119, 7).Value = "text"
$ws.Cells.Item(119, 8).Value = "DFDL-199"

$ws.Cells.Item(120, 2).Value = "text_05"
$ws.Cells.Item(120, 3).Value = "Passed"
$ws.Cells.Item(120, 4).Value = "DFDL-13-235R"
$ws.Cells.Item(120, 5).Value = "High"
$ws.Cells.Item(120, 6).Value = 3
$ws.Cells.Item(120, 7).Value = "text"
$ws.Cells.Item(120, 8).Value = "DFDL-199"

$ws.Cells.Item(121, 2).Value = "text_06"
$ws.Cells.Item(121, 3).Value = "Passed"
$ws.Cells.Item(121, 4).Value = "DFDL-13-235R"
$ws.Cells.Item(121, 5).Value = "High"
$ws.Cells.Item(121, 6).Value = 3
$ws.Cells.Item(121, 7).Value = "text"
$ws.Cells.Item(121, 8).Value = "DFDL-199"

$ws.Cells.Item(122, 2).Value = "binary_01"
$ws.Cells.Item(122, 3).Value = "Passed"
$ws.Cells.Item(122, 4).Value = "DFDL-13-235R"
$ws.Cells.Item(122, 5).Value = "High"
$ws.Cells.Item(122, 6).Value = 3
$ws.Cells.Item(122, 7).Value = "byte"
$ws.Cells.Item(122, 8).Value = "DFDL-199"

$ws.Cells.Item(123, 2).Value = "entity_fail_05"
$ws.Cells.Item(123, 3).Value = "Passed"
$ws.Cells.Item(123, 4).Value = "DFDL-13-235R"
$ws.Cells.Item(123, 5).Value = "High"
$ws.Cells.Item(123, 6).Value = 3
$ws.Cells.Item(123, 7).Value = "byte"
$ws.Cells.Item(123, 8).Value = "DFDL-199"
$ws.Cells.Item(123, 9).Value = "Yes"

$ws.Cells.Item(124, 2).Value = "entity_fail_06"
$ws.Cells.Item(124, 3).Value = "Passed"
$ws.Cells.Item(124, 4).Value = "DFDL-13-235R"
$ws.Cells.Item(124, 5).Value = "High"
$ws.Cells.Item(124, 6).Value = 3
$ws.Cells.Item(124, 7).Value = "byte"
$ws.Cells.Item(124, 8).Value = "DFDL-199"
$ws.Cells.Item(124, 9).Value = "Yes"

$ws.Cells.Item(125, 2).Value = "property_scoping_06"
$ws.Cells.Item(125, 3).Value = "Passed"
$ws.Cells.Item(125, 4).Value = "DFDL-8-022R"
$ws.Cells.Item(125, 5).Value = "High"
$ws.Cells.Item(125, 6).Value = 4
$ws.Cells.Item(125, 7).Value = "text"
$ws.Cells.Item(125, 8).Value = "DFDL-281"

$ws.Cells.Item(126, 2).Value = "NumSeq_05"
$ws.Cells.Item(126, 3).Value = "Passed"
$ws.Cells.Item(126, 4).Value = "DFDL-12-043R"
$ws.Cells.Item(126, 5).Value = "High"
$ws.Cells.Item(126, 6).Value = 4
$ws.Cells.Item(126, 7).Value = "text"
$ws.Cells.Item(126, 8).Value = "DFDL-63"

$ws.Cells.Item(127, 2).Value = "NumSeq_06"
$ws.Cells.Item(127, 3).Value = "Passed"
$ws.Cells.Item(127, 4).Value = "DFDL-12-043R"
$ws.Cells.Item(127, 5).Value = "High"
$ws.Cells.Item(127, 6).Value = 4
$ws.Cells.Item(127, 7).Value = "text"
$ws.Cells.Item(127, 8).Value = "DFDL-63"

$ws.Cells.Item(128, 2).Value = "NumSeq_07"
$ws.Cells.Item(128, 3).Value = "Passed"
$ws.Cells.Item(128, 4).Value = "DFDL-12-043R"
$ws.Cells.Item(128, 5).Value = "High"
$ws.Cells.Item(128, 6).Value = 4
$ws.Cells.Item(128, 7).Value = "text"
$ws.Cells.Item(128, 8).Value = "DFDL-63"

$ws.Cells.Item(129, 2).Value = "NumSeq_08"
$ws.Cells.Item(129, 3).Value = "Passed"
$ws.Cells.Item(129, 4).Value = "DFDL-12-043R"
$ws.Cells.Item(129, 5).Value = "High"
$ws.Cells.Item(129, 6).Value = 4
$ws.Cells.Item(129, 7).Value = "text"
$ws.Cells.Item(129, 8).Value = "DFDL-63"

# --- Update the view: selection moves to F129, scrolled so B103 is the top-left visible cell ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 103
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F129").Select()

# --- Restore workbook window placement/size (best effort) ---
$excel.Windows.Item(1).Top = 636
$excel.Windows.Item(1).Height = 11280
